$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "CreatedAt: 2025-05-07T18:07:14"
$ws.Range("V4").Value = 183.43
$ws.Range("W4").Value = 171.65
$ws.Range("X4").Value = 14.64
$ws.Range("Y4").Value = 6.62
$ws.Range("Z4").Value = 6.64
$ws.Range("V6").Value = 1.1
$ws.Range("W6").Value = 1.37
$ws.Range("V9").Value = 183.98
$ws.Range("W9").Value = 167.76
$ws.Range("X9").Value = 15.07
$ws.Range("Y9").Value = 6.78
$ws.Range("Z9").Value = 6.76
$ws.Range("V11").Value = 1.66
$ws.Range("W11").Value = -2.52
$ws.Range("X11").Value = 0.59
$ws.Range("Y11").Value = 0.28
$ws.Range("V14").Value = 183.98
$ws.Range("W14").Value = 167.59
$ws.Range("X14").Value = 15.07
$ws.Range("Y14").Value = 6.78
$ws.Range("Z14").Value = 6.77
$ws.Range("V16").Value = 1.66
$ws.Range("W16").Value = -2.68
$ws.Range("X16").Value = 0.59
$ws.Range("Y16").Value = 0.28
$ws.Range("V19").Value = 184.91
$ws.Range("W19").Value = 173.04
$ws.Range("X19").Value = 14.76
$ws.Range("Y19").Value = 6.66
$ws.Range("Z19").Value = 6.69
$ws.Range("V21").Value = 2.59
$ws.Range("W21").Value = 2.77
$ws.Range("V24").Value = 184.91
$ws.Range("W24").Value = 173.04
$ws.Range("X24").Value = 14.76
$ws.Range("Y24").Value = 6.66
$ws.Range("Z24").Value = 6.69
$ws.Range("V26").Value = 2.59
$ws.Range("W26").Value = 2.77
$ws.Range("V29").Value = 187
$ws.Range("W29").Value = 175
$ws.Range("X29").Value = 14.9
$ws.Range("Y29").Value = 6.7
$ws.Range("Z29").Value = 6.72
$ws.Range("V31").Value = 4.68
$ws.Range("W31").Value = 4.73
$ws.Range("Y31").Value = 0.2
$ws.Range("Z31").Value = 0.23
$ws.Range("V34").Value = 181.6
$ws.Range("W34").Value = 165.16
$ws.Range("X34").Value = 15.08
$ws.Range("Y34").Value = 6.78
$ws.Range("Z34").Value = 6.73
$ws.Range("V36").Value = -0.73
$ws.Range("W36").Value = -5.12
$ws.Range("Y36").Value = 0.28
$ws.Range("V39").Value = 183.43
$ws.Range("W39").Value = 171.65
$ws.Range("X39").Value = 14.64
$ws.Range("Y39").Value = 6.62
$ws.Range("Z39").Value = 6.64
$ws.Range("V41").Value = 1.1
$ws.Range("W41").Value = 1.37
$ws.Range("V44").Value = 183.98
$ws.Range("W44").Value = 171.99
$ws.Range("X44").Value = 14.58
$ws.Range("Y44").Value = 6.56
$ws.Range("Z44").Value = 6.56
$ws.Range("V46").Value = 1.66
$ws.Range("W46").Value = 1.72
$ws.Range("V49").Value = 168.51
$ws.Range("W49").Value = 156.79
$ws.Range("X49").Value = 13.27
$ws.Range("Y49").Value = 5.92
$ws.Range("Z49").Value = 6.19
$ws.Range("V51").Value = -13.82
$ws.Range("W51").Value = -13.48
$ws.Range("Y51").Value = -0.57
$ws.Range("V54").Value = 162.21
$ws.Range("W54").Value = 151.76
$ws.Range("X54").Value = 13.44
$ws.Range("Y54").Value = 6.02
$ws.Range("Z54").Value = 5.97
$ws.Range("V55").Value = 0
$ws.Range("V56").Value = -20.11
$ws.Range("W56").Value = -18.51
$ws.Range("Y56").Value = -0.48
$ws.Range("V59").Value = 189.53
$ws.Range("W59").Value = 177
$ws.Range("X59").Value = 14.96
$ws.Range("Y59").Value = 6.73
$ws.Range("Z59").Value = 6.72
$ws.Range("V61").Value = 7.2
$ws.Range("W61").Value = 6.73
$ws.Range("X61").Value = 0.48
$ws.Range("Y61").Value = 0.23
$ws.Range("V64").Value = 192.33
$ws.Range("W64").Value = 179.61
$ws.Range("X64").Value = 15.16
$ws.Range("Y64").Value = 6.8
$ws.Range("Z64").Value = 6.79
$ws.Range("V66").Value = 10
$ws.Range("W66").Value = 9.34
$ws.Range("Y66").Value = 0.31
$ws.Range("V69").Value = 194.38
$ws.Range("W69").Value = 181.72
$ws.Range("X69").Value = 15.31
$ws.Range("Y69").Value = 6.86
$ws.Range("Z69").Value = 6.85
$ws.Range("V71").Value = 12.05
$ws.Range("W71").Value = 11.45
$ws.Range("X71").Value = 0.83
$ws.Range("Y71").Value = 0.36
$ws.Range("V74").Value = 173.81
$ws.Range("W74").Value = 162.32
$ws.Range("X74").Value = 13.88
$ws.Range("Y74").Value = 6.22
$ws.Range("Z74").Value = 6.25
$ws.Range("V76").Value = -8.52
$ws.Range("W76").Value = -7.95
$ws.Range("Y76").Value = -0.27
$ws.Range("V79").Value = 182.32
$ws.Range("W79").Value = 170.27
$ws.Range("X79").Value = 14.48
$ws.Range("Y79").Value = 6.5
$ws.Range("Z79").Value = 6.49
$ws.Range("V84").Value = 152.96
$ws.Range("W84").Value = 142.97
$ws.Range("X84").Value = 13.62
$ws.Range("Y84").Value = 6.1
$ws.Range("Z84").Value = 6.07
$ws.Range("V85").Value = 0
$ws.Range("V86").Value = -29.37
$ws.Range("W86").Value = -27.31
$ws.Range("Y86").Value = -0.4
$ws.Range("Z86").Value = -0.42
$ws.Range("V89").Value = 187
$ws.Range("W89").Value = 175
$ws.Range("X89").Value = 14.9
$ws.Range("Y89").Value = 6.7
$ws.Range("Z89").Value = 6.72
$ws.Range("V91").Value = 4.68
$ws.Range("W91").Value = 4.73
$ws.Range("Y91").Value = 0.2
$ws.Range("Z91").Value = 0.23

Write-Host "Applied all changes"